$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17, shifting existing rows 17-27 down to 18-28
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with data (same as the other Perejil rows, but with
# a new date and "Primera" quality figures)
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = (Get-Date -Year 2022 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 100112044
$ws.Cells.Item(17, 7).Value = "Perejil"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 200
$ws.Cells.Item(17, 11).Value = 700
$ws.Cells.Item(17, 12).Value = 800
$ws.Cells.Item(17, 13).Value = 750
$ws.Cells.Item(17, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(17, 15).Value = "Región del Maule"
$ws.Cells.Item(17, 16).Value = 750
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"
